$wb = $excel.ActiveWorkbook

$dateValue = [double]"45835.49694444444"

# --- Sheet 1: FE_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$r = 49
$ws.Cells.Item($r,1).Value = $dateValue
$ws.Cells.Item($r,1).NumberFormat = $ws.Cells.Item($r-1,1).NumberFormat
$ws.Cells.Item($r,2).Value = "0x01,0x7c"
$ws.Cells.Item($r,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item($r,4).Value = "0x01,0x5C"
$ws.Cells.Item($r,5).Value = "0xf"
$ws.Cells.Item($r,6).Value = 380
$ws.Cells.Item($r,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item($r,8).Value = 348
$ws.Cells.Item($r,9).Value = 15

# --- Sheet 2: FE_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$r = 49
$ws.Cells.Item($r,1).Value = $dateValue
$ws.Cells.Item($r,1).NumberFormat = $ws.Cells.Item($r-1,1).NumberFormat
$ws.Cells.Item($r,2).Value = "0x01,0x90"
$ws.Cells.Item($r,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item($r,4).Value = "0x01,0x6C"
$ws.Cells.Item($r,5).Value = "0xe"
$ws.Cells.Item($r,6).Value = 400
$ws.Cells.Item($r,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item($r,8).Value = 364
$ws.Cells.Item($r,9).Value = 14

# --- Sheet 3: FE_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$r = 49
$ws.Cells.Item($r,1).Value = $dateValue
$ws.Cells.Item($r,1).NumberFormat = $ws.Cells.Item($r-1,1).NumberFormat
$ws.Cells.Item($r,2).Value = "0x00,0x6e"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x69"
$ws.Cells.Item($r,5).Value = "0x3"
$ws.Cells.Item($r,6).Value = 110
$ws.Cells.Item($r,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r,8).Value = 105
$ws.Cells.Item($r,9).Value = 3

# --- Sheet 4: FE_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$r = 49
$ws.Cells.Item($r,1).Value = $dateValue
$ws.Cells.Item($r,1).NumberFormat = $ws.Cells.Item($r-1,1).NumberFormat
$ws.Cells.Item($r,2).Value = "0x00,0x6e"
$ws.Cells.Item($r,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x69"
$ws.Cells.Item($r,5).Value = "0x3"
$ws.Cells.Item($r,6).Value = 110
$ws.Cells.Item($r,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item($r,8).Value = 105
$ws.Cells.Item($r,9).Value = 3
